# DepartmentMeetings_template.xlsx edit
# 1. Add a new "number of attendees" column (L) to the template:
#    - L1 header cell: "参会人数" with the same header style as the
#      existing header cells (copy format from K1, the last header cell).
#    - L2 data cell: the merge-field placeholder "${record.numberOfPeople}".
# 2. Give the new column a sensible width.
# 3. Move the active selection (this is simply where the author's cursor
#    happened to be left when the workbook was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the header style from K1 onto the new L1 header cell ---------
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)   # xlPasteFormats

# --- Header text for the new column ------------------------------------
$ws.Cells.Item(1, 12).Value = "参会人数"

# --- Merge-field placeholder for the new column's data row -------------
$ws.Cells.Item(2, 12).Value = "`${record.numberOfPeople}"

# --- Give column L (12) an explicit width -------------------------------
$ws.Columns("L").ColumnWidth = 8.92

# --- Restore the cursor/selection left behind in the saved file --------
$ws.Range("F12").Select() | Out-Null
